$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'306.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.65%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'39.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'7.45%"
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'0.76%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08051"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-0.82%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.955"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-5.51%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'4.196"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'1.01%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'8.001"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'1.92%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9324"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.39%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1448"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'2.44%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1931"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-0.20%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.09150"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'0.11%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03509"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'1.42%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09792"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-1.35%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001401"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.27%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.006034"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-4.48%"
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'-1.36%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'2.25%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3423"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-0.62%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1304"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-0.24%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.558"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-5.29%"
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'3.25%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04376"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.04%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'0.34%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004272"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-13.15%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001302"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'0.18%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D39").Value = "'0.02026"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'0.05%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.05053"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-2.17%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007445"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-0.79%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.01030"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'1.45%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.1345"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-1.95%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002124"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-2.13%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.009111"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-8.65%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006196"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-1.27%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'0.19%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.003098"
$ws.Range("D48").Style = "Normal"
$ws.Range("E49").Value = "'28.27%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'0.19%"
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'0.19%"
$ws.Range("E51").Style = "Normal"
